$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics for Dhh-Ptch1 ligand-receptor pairs
# (recomputed with new TPM values; only the affected numeric cells are set)

# Row 2
$ws.Range("G2").Value = 3.659958333333333
$ws.Range("H2").Value = 10.979875
$ws.Range("I2").Value = 0.4781132044744068
$ws.Range("J2").Value = 0.4781132044744067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.695620666666667
$ws.Range("N2").Value = 23.086862
$ws.Range("O2").Value = 0.4976976897997125
$ws.Range("P2").Value = 0.4976976897997126
$ws.Range("Q2").Value = 28.16565098913889
$ws.Range("R2").Value = 253.49085890225
$ws.Range("S2").Value = 0.2379558373296498
$ws.Range("T2").Value = 0.2379558373296498

# Row 3
$ws.Range("G3").Value = 3.659958333333333
$ws.Range("H3").Value = 10.979875
$ws.Range("I3").Value = 0.4781132044744068
$ws.Range("J3").Value = 0.4781132044744067
$ws.Range("O3").Value = 0.228990810419744
$ws.Range("P3").Value = 0.228990810419744
$ws.Range("Q3").Value = 12.95902186847222
$ws.Range("R3").Value = 116.63119681625
$ws.Range("S3").Value = 0.1094835301649752
$ws.Range("T3").Value = 0.1094835301649752

# Row 4
$ws.Range("G4").Value = 3.659958333333333
$ws.Range("H4").Value = 10.979875
$ws.Range("I4").Value = 0.4781132044744068
$ws.Range("J4").Value = 0.4781132044744067
$ws.Range("M4").Value = 3.610968333333334
$ws.Range("N4").Value = 10.832905
$ws.Range("O4").Value = 0.2335315987213747
$ws.Range("P4").Value = 0.2335315987213747
$ws.Range("Q4").Value = 13.21599364298611
$ws.Range("R4").Value = 118.943942786875
$ws.Range("S4").Value = 0.1116545410107078
$ws.Range("T4").Value = 0.1116545410107077

# Row 5
$ws.Range("G5").Value = 3.659958333333333
$ws.Range("H5").Value = 10.979875
$ws.Range("I5").Value = 0.4781132044744068
$ws.Range("J5").Value = 0.4781132044744067
$ws.Range("M5").Value = 0.6150943333333334
$ws.Range("N5").Value = 1.845283
$ws.Range("O5").Value = 0.03977990105916876
$ws.Range("P5").Value = 0.03977990105916877
$ws.Range("Q5").Value = 2.251219631069445
$ws.Range("R5").Value = 20.260976679625
$ws.Range("S5").Value = 0.01901929596907402
$ws.Range("T5").Value = 0.01901929596907402

# Row 6
$ws.Range("I6").Value = 0.3193330932870009
$ws.Range("J6").Value = 0.3193330932870008
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.695620666666667
$ws.Range("N6").Value = 23.086862
$ws.Range("O6").Value = 0.4976976897997125
$ws.Range("P6").Value = 0.4976976897997126
$ws.Range("Q6").Value = 18.811913937184
$ws.Range("R6").Value = 169.307225434656
$ws.Range("S6").Value = 0.1589313428055364
$ws.Range("T6").Value = 0.1589313428055364

# Row 7
$ws.Range("I7").Value = 0.3193330932870009
$ws.Range("J7").Value = 0.3193330932870008
$ws.Range("O7").Value = 0.228990810419744
$ws.Range("P7").Value = 0.228990810419744
$ws.Range("S7").Value = 0.07312434382563404
$ws.Range("T7").Value = 0.07312434382563404

# Row 8
$ws.Range("I8").Value = 0.3193330932870009
$ws.Range("J8").Value = 0.3193330932870008
$ws.Range("M8").Value = 3.610968333333334
$ws.Range("N8").Value = 10.832905
$ws.Range("O8").Value = 0.2335315987213747
$ws.Range("P8").Value = 0.2335315987213747
$ws.Range("Q8").Value = 8.826997646960001
$ws.Range("R8").Value = 79.44297882264
$ws.Range("S8").Value = 0.0745743677999552
$ws.Range("T8").Value = 0.0745743677999552

# Row 9
$ws.Range("I9").Value = 0.3193330932870009
$ws.Range("J9").Value = 0.3193330932870008
$ws.Range("M9").Value = 0.6150943333333334
$ws.Range("N9").Value = 1.845283
$ws.Range("O9").Value = 0.03977990105916876
$ws.Range("P9").Value = 0.03977990105916877
$ws.Range("Q9").Value = 1.503595637456
$ws.Range("R9").Value = 13.532360737104
$ws.Range("S9").Value = 0.0127030388558752
$ws.Range("T9").Value = 0.0127030388558752

# Row 10
$ws.Range("G10").Value = 1.388093333333333
$ws.Range("H10").Value = 4.16428
$ws.Range("I10").Value = 0.1813315046964271
$ws.Range("J10").Value = 0.1813315046964271
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.695620666666667
$ws.Range("N10").Value = 23.086862
$ws.Range("O10").Value = 0.4976976897997125
$ws.Range("P10").Value = 0.4976976897997126
$ws.Range("Q10").Value = 10.68223974326222
$ws.Range("R10").Value = 96.14015768936
$ws.Range("S10").Value = 0.09024827097531748
$ws.Range("T10").Value = 0.0902482709753175

# Row 11
$ws.Range("G11").Value = 1.388093333333333
$ws.Range("H11").Value = 4.16428
$ws.Range("I11").Value = 0.1813315046964271
$ws.Range("J11").Value = 0.1813315046964271
$ws.Range("O11").Value = 0.228990810419744
$ws.Range("P11").Value = 0.228990810419744
$ws.Range("Q11").Value = 4.914900723955555
$ws.Range("R11").Value = 44.2341065156
$ws.Range("S11").Value = 0.04152324821506646
$ws.Range("T11").Value = 0.04152324821506646

# Row 12
$ws.Range("G12").Value = 1.388093333333333
$ws.Range("H12").Value = 4.16428
$ws.Range("I12").Value = 0.1813315046964271
$ws.Range("J12").Value = 0.1813315046964271
$ws.Range("M12").Value = 3.610968333333334
$ws.Range("N12").Value = 10.832905
$ws.Range("O12").Value = 0.2335315987213747
$ws.Range("P12").Value = 0.2335315987213747
$ws.Range("Q12").Value = 5.012361070377778
$ws.Range("R12").Value = 45.1112496334
$ws.Range("S12").Value = 0.04234663619030909
$ws.Range("T12").Value = 0.04234663619030909

# Row 13
$ws.Range("G13").Value = 1.388093333333333
$ws.Range("H13").Value = 4.16428
$ws.Range("I13").Value = 0.1813315046964271
$ws.Range("J13").Value = 0.1813315046964271
$ws.Range("M13").Value = 0.6150943333333334
$ws.Range("N13").Value = 1.845283
$ws.Range("O13").Value = 0.03977990105916876
$ws.Range("P13").Value = 0.03977990105916877
$ws.Range("Q13").Value = 0.8538083434711111
$ws.Range("R13").Value = 7.684275091240001
$ws.Range("S13").Value = 0.007213349315734066
$ws.Range("T13").Value = 0.007213349315734066

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.162456
$ws.Range("H14").Value = 0.487368
$ws.Range("I14").Value = 0.02122219754216535
$ws.Range("J14").Value = 0.02122219754216534
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 7.695620666666667
$ws.Range("N14").Value = 23.086862
$ws.Range("O14").Value = 0.4976976897997125
$ws.Range("P14").Value = 0.4976976897997126
$ws.Range("Q14").Value = 1.250199751024
$ws.Range("R14").Value = 11.251797759216
$ws.Range("S14").Value = 0.01056223868920883
$ws.Range("T14").Value = 0.01056223868920883

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.162456
$ws.Range("H15").Value = 0.487368
$ws.Range("I15").Value = 0.02122219754216535
$ws.Range("J15").Value = 0.02122219754216534
$ws.Range("O15").Value = 0.228990810419744
$ws.Range("P15").Value = 0.228990810419744
$ws.Range("Q15").Value = 0.5752171650400001
$ws.Range("R15").Value = 5.17695448536
$ws.Range("S15").Value = 0.004859688214068342
$ws.Range("T15").Value = 0.004859688214068342

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.162456
$ws.Range("H16").Value = 0.487368
$ws.Range("I16").Value = 0.02122219754216535
$ws.Range("J16").Value = 0.02122219754216534
$ws.Range("M16").Value = 3.610968333333334
$ws.Range("N16").Value = 10.832905
$ws.Range("O16").Value = 0.2335315987213747
$ws.Range("P16").Value = 0.2335315987213747
$ws.Range("Q16").Value = 0.5866234715600001
$ws.Range("R16").Value = 5.279611244040001
$ws.Range("S16").Value = 0.004956053720402702
$ws.Range("T16").Value = 0.004956053720402702

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.162456
$ws.Range("H17").Value = 0.487368
$ws.Range("I17").Value = 0.02122219754216535
$ws.Range("J17").Value = 0.02122219754216534
$ws.Range("M17").Value = 0.6150943333333334
$ws.Range("N17").Value = 1.845283
$ws.Range("O17").Value = 0.03977990105916876
$ws.Range("P17").Value = 0.03977990105916877
$ws.Range("Q17").Value = 0.09992576501600002
$ws.Range("R17").Value = 0.8993318851440002
$ws.Range("S17").Value = 0.0008442169184854719
$ws.Range("T17").Value = 0.0008442169184854719

